# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "29.806.64"
Set-TextValue "E2" "  -1.35%  "

Set-TextValue "D3" "1.890.44"
Set-TextValue "E3" "  -0.99%  "

Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.22%  "

Set-TextValue "D5" "0.7534"
Set-TextValue "E5" "  +2.98%  "

Set-TextValue "D6" "239.58"
Set-TextValue "E6" "  -1.79%  "

Set-TextValue "E7" "  -0.12%  "

Set-TextValue "D8" "0.3034"
Set-TextValue "E8" "  -3.14%  "

Set-TextValue "D9" "25.43"
Set-TextValue "E9" "  -5.10%  "

Set-TextValue "D10" "0.06803"

Set-TextValue "D11" "0.07946"
Set-TextValue "E11" "  -0.49%  "

Set-TextValue "B12" "Polygon"
Set-TextValue "C12" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D12" "0.7412"
Set-TextValue "E12" "  -4.68%  "

Set-TextValue "B13" "WrappedEther"
Set-TextValue "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.891.60"
Set-TextValue "E13" "  +0.30%  "

Set-TextValue "D14" "5.161"
Set-TextValue "E14" "  -1.77%  "

Set-TextValue "D15" "90.70"
Set-TextValue "E15" "  -0.78%  "

Set-TextValue "D16" "29.807.11"
Set-TextValue "E16" "  -1.23%  "

Set-TextValue "E17" "  -2.49%  "

Set-TextValue "D18" "5.938"
Set-TextValue "E18" "  +1.61%  "

Set-TextValue "D19" "243.90"
Set-TextValue "E19" "  +1.67%  "

Set-TextValue "D20" "0.000007669"
Set-TextValue "E20" "  -1.80%  "

Set-TextValue "E21" "  -0.26%  "

Set-TextValue "D22" "1.001"
Set-TextValue "E22" "  -0.23%  "

Set-TextValue "D23" "6.911"
Set-TextValue "E23" "  +1.88%  "

Set-TextValue "D24" "9.211"
Set-TextValue "E24" "  -1.91%  "

Set-TextValue "D25" "165.37"
Set-TextValue "E25" "  -0.22%  "

Set-TextValue "D26" "18.65"
Set-TextValue "E26" "  -2.41%  "

Set-TextValue "D27" "0.1272"
Set-TextValue "E27" "  +0.21%  "

Set-TextValue "D28" "2.016"
Set-TextValue "E28" "  -3.49%  "

Set-TextValue "D29" "1.380"
Set-TextValue "E29" "  +1.98%  "

Set-TextValue "D30" "1.515"
Set-TextValue "E30" "  -2.17%  "

Set-TextValue "D31" "4.234"
Set-TextValue "E31" "  -1.64%  "

Set-TextValue "D32" "3.993"
Set-TextValue "E32" "  -2.28%  "

Set-TextValue "D33" "0.05250"
Set-TextValue "E33" "  +1.76%  "

Set-TextValue "D34" "1.242"
Set-TextValue "E34" "  -3.19%  "

Set-TextValue "D35" "0.7191"
Set-TextValue "E35" "  -3.14%  "

Set-TextValue "D36" "2.717"

Set-TextValue "D37" "0.01899"
Set-TextValue "E37" "  -2.00%  "

Set-TextValue "D38" "2.771"
Set-TextValue "E38" "  -0.93%  "

Set-TextValue "D39" "6.137"
Set-TextValue "E39" "  -3.31%  "

Set-TextValue "D40" "0.4380"
Set-TextValue "E40" "  -1.31%  "

Set-TextValue "D41" "71.64"
Set-TextValue "E41" "  -4.20%  "

Set-TextValue "D42" "1.002"
Set-TextValue "E42" "  -0.04%  "

Set-TextValue "D43" "1.876"
Set-TextValue "E43" "  -2.82%  "

Set-TextValue "D44" "0.8246"
Set-TextValue "E44" "  -1.29%  "

Set-TextValue "D45" "100.47"
Set-TextValue "E45" "  -0.62%  "

Set-TextValue "B46" "Aptos"
Set-TextValue "C46" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D46" "7.504"
Set-TextValue "E46" "  -1.24%  "

Set-TextValue "B47" "EnergySwap"
Set-TextValue "C47" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "9.717"
Set-TextValue "E47" "  -0.28%  "

Set-TextValue "D48" "2.049.78"
Set-TextValue "E48" "  -0.03%  "

Set-TextValue "D49" "35.91"
Set-TextValue "E49" "  -4.16%  "

Set-TextValue "E50" "  -0.40%  "

Set-TextValue "D51" "1.462"
Set-TextValue "E51" "  +0.39%  "
